# Edit script: "cleaned version of GA"
# Applies content, formatting and merge-cell changes to the
# "Group 2 Preparatory" timetable sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Color constants (OLE/BGR integers understood by Range.Interior.Color)
# ---------------------------------------------------------------------
$GREEN  = 5296274   # 92D050
$YELLOW = 65535      # FFFF00
$WHITE  = 16777215   # FFFFFF
$RED    = 255        # FF0000

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Simple text-only edit: F6 room number change
# ---------------------------------------------------------------------
$ws.Range("F6").Value = "Academic Writing: Research, Fiction and Nonfiction`n`n09:00-10:45`nroom:209"
$ws.Range("F6").Interior.Color = $YELLOW

# ---------------------------------------------------------------------
# 2. Row 14-21 block (Monday/Wed/Thu/Fri columns B/D/E/F)
#    Old merges: B14:B20, D14:D20, F14:F19
#    New merges: B14:B21, D14:D21, E14:E20, F14:F21
# ---------------------------------------------------------------------
$ws.Range("B14:B20").UnMerge()
$ws.Range("D14:D20").UnMerge()
$ws.Range("F14:F19").UnMerge()

$ws.Range("B14:B21").Merge()
$ws.Range("D14:D21").Merge()
$ws.Range("E14:E20").Merge()
$ws.Range("F14:F21").Merge()

$ws.Range("B14").Value = "Scientific Inquiry: Beyond the Visible`n`n11:00-12:30`nroom:204"
$ws.Range("B14").Interior.Color = $GREEN

$ws.Range("D14").Value = "Precalculus`n`n11:00-13:00`nroom:204"
$ws.Range("D14").Interior.Color = $GREEN

$ws.Range("E14").Value = "Scientific Inquiry: Beyond the Visible`n`n11:00-12:30`nroom:104"
$ws.Range("E14").Interior.Color = $WHITE

$ws.Range("F14").Value = "Precalculus`n`n11:00-13:00`nroom:202"
$ws.Range("F14").Interior.Color = $WHITE

# Borders / fill for the empty filler cells inside the merged blocks.
# B19 / E19 become "bottom border" style (copy from former B20 pattern),
# while B20 / D20 / F19 become plain "no-bottom-border" middle cells,
# and row 21 (new bottom row of each merge) gets the bottom-border style.
$ws.Range("D19").Copy() | Out-Null
$ws.Range("B19").PasteSpecial($xlPasteFormats)
$ws.Range("E19").PasteSpecial($xlPasteFormats)

$ws.Range("B15").Copy() | Out-Null
$ws.Range("B20").PasteSpecial($xlPasteFormats)
$ws.Range("D20").PasteSpecial($xlPasteFormats)
$ws.Range("F19").PasteSpecial($xlPasteFormats)
$ws.Range("F20").PasteSpecial($xlPasteFormats)
$ws.Range("E15").PasteSpecial($xlPasteFormats)
$ws.Range("E16").PasteSpecial($xlPasteFormats)
$ws.Range("E17").PasteSpecial($xlPasteFormats)
$ws.Range("E18").PasteSpecial($xlPasteFormats)

$ws.Range("B20").Copy() | Out-Null
$ws.Range("E20").PasteSpecial($xlPasteFormats)

$ws.Range("D20").Copy() | Out-Null
$ws.Range("B21").PasteSpecial($xlPasteFormats)
$ws.Range("D21").PasteSpecial($xlPasteFormats)
$ws.Range("F21").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Row 26 block
#    Old merges: F26:F32 (keep/extend), D26:D33 (removed), B26:B32 (kept)
# ---------------------------------------------------------------------
$ws.Range("D26:D33").UnMerge()
$ws.Range("F26:F32").UnMerge()
$ws.Range("F26:F37").Merge()

$ws.Range("B26").Value = "Academic Writing: Research, Fiction and Nonfiction`n`n14:00-15:45`nroom:204"
$ws.Range("B26").Interior.Color = $GREEN

$ws.Range("D26").Clear() | Out-Null

$ws.Range("F26").Value = "Scientific Inquiry: Beyond the Visible`n`n14:00-17:00`nroom:-14"
$ws.Range("F26").Interior.Color = $WHITE

# Clear stray D cells (27-33) left over from the removed D26:D33 merge
foreach ($r in 27..33) {
    $ws.Range("D$r").Clear() | Out-Null
}

# F32 loses its old bottom-border look (no longer end of merge); F33 becomes
# a plain filler cell inside the extended F26:F37 merge.
$ws.Range("F27").Copy() | Out-Null
$ws.Range("F32").PasteSpecial($xlPasteFormats)
$ws.Range("F33").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. Row 34-40 block
#    Old merges: B34:B40 (kept), D34:D39->D34:D40, F34:F40 (removed)
#    New merges: C34:C40, E34:E40 (added)
# ---------------------------------------------------------------------
$ws.Range("D34:D39").UnMerge()
$ws.Range("F34:F40").UnMerge()

$ws.Range("D34:D40").Merge()
$ws.Range("C34:C40").Merge()
$ws.Range("E34:E40").Merge()

$ws.Range("B34").Value = "Academic Writing: Research, Fiction and Nonfiction`n`n16:00-17:45`nroom:204"
$ws.Range("B34").Interior.Color = $GREEN

$ws.Range("C34").Value = "Scientific Inquiry: Beyond the Visible`n`n16:00-17:30`nroom:202"
$ws.Range("C34").Interior.Color = $WHITE

$ws.Range("D34").Value = "Academic Writing: Research, Fiction and Nonfiction`n`n16:00-17:45`nroom:209"
$ws.Range("D34").Interior.Color = $YELLOW

$ws.Range("E34").Value = "Academic Writing: Research, Fiction and Nonfiction`n`n16:00-17:45`nroom:209"
$ws.Range("E34").Interior.Color = $YELLOW

$ws.Range("F34").Clear() | Out-Null
$ws.Range("F35").Copy() | Out-Null
$ws.Range("F34").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# Fill the new C/E columns (rows 35-40) with plain filler style, matching
# the existing B/D pattern for that row.
foreach ($r in 35..40) {
    $ws.Range("B$r").Copy() | Out-Null
    $ws.Range("C$r").PasteSpecial($xlPasteFormats)
    $ws.Range("E$r").PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false
}

# Row-specific border tweaks: C39 and F37 get the bottom-border look,
# F38/F39/F40 are no longer part of any merge (cleared), D40 gets bottom
# border (end of the D34:D40 merge), C40/E40/B40 already end in bottom
# border style (B40 unchanged, copy for C40/E40 done above already needs fix)
$ws.Range("B39").Copy() | Out-Null
$ws.Range("C39").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("B40").Copy() | Out-Null
$ws.Range("C40").PasteSpecial($xlPasteFormats)
$ws.Range("D40").PasteSpecial($xlPasteFormats)
$ws.Range("E40").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("B32").Copy() | Out-Null
$ws.Range("F37").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

foreach ($r in 38..40) {
    $ws.Range("F$r").Clear() | Out-Null
}

# ---------------------------------------------------------------------
# 5. Rows 41-45: stray D cells removed (D34:D45 shrunk to D34:D40)
# ---------------------------------------------------------------------
foreach ($r in 41..45) {
    $ws.Range("D$r").Clear() | Out-Null
}
